# Apply Ifrit_Profits.xlsx market-price refresh values (scheduled runner update).
# Values below come from the authoritative diff; columns H-N are derived
# market-price / profit figures (currentAveragePrice*, LevePrice*, LeveProfit*).
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 549.5
$ws.Range("I6").Value = 100
$ws.Range("J6").Value = 819.2
$ws.Range("K6").Value = 300
$ws.Range("L6").Value = 2457.6
$ws.Range("M6").Value = -188
$ws.Range("N6").Value = -2681.6
$ws.Range("H62").Value = 2435.6445
$ws.Range("I62").Value = 1074.1428
$ws.Range("K62").Value = 1074.1428
$ws.Range("M62").Value = -450.1428000000001
$ws.Range("H65").Value = 2435.6445
$ws.Range("I65").Value = 1074.1428
$ws.Range("K65").Value = 5370.714
$ws.Range("M65").Value = -2250.714
$ws.Range("H70").Value = 253250.75
$ws.Range("H73").Value = 253250.75
$ws.Range("H76").Value = 3203.125
$ws.Range("J76").Value = 3003.5
$ws.Range("L76").Value = 3003.5
$ws.Range("N76").Value = -3633.5
$ws.Range("H79").Value = 3203.125
$ws.Range("J79").Value = 3003.5
$ws.Range("L79").Value = 3003.5
$ws.Range("N79").Value = -5187.5
$ws.Range("H131").Value = 1687.9166
$ws.Range("I131").Value = 1687.9166
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 5063.7498
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = -23.7497999999996
$ws.Range("H132").Value = 1050.3334
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 1050.3334
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 3151.0002
$ws.Range("N132").Value = -8211.0002
$ws.Range("H138").Value = 2741.1628
$ws.Range("I138").Value = 2718.6667
$ws.Range("J138").Value = 2757.36
$ws.Range("K138").Value = 8156.000100000001
$ws.Range("L138").Value = 8272.08
$ws.Range("M138").Value = -3016.000100000001
$ws.Range("N138").Value = -18552.08
$ws.Range("N131").ClearContents()
$ws.Range("M132").ClearContents()

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2682.4546
$ws.Range("I88").Value = 2789
$ws.Range("J88").Value = 2398.3333
$ws.Range("K88").Value = 2789
$ws.Range("L88").Value = 2398.3333
$ws.Range("M88").Value = -2383
$ws.Range("N88").Value = -3210.3333
$ws.Range("H91").Value = 2682.4546
$ws.Range("I91").Value = 2789
$ws.Range("J91").Value = 2398.3333
$ws.Range("K91").Value = 2789
$ws.Range("L91").Value = 2398.3333
$ws.Range("M91").Value = -1385
$ws.Range("N91").Value = -5206.3333
$ws.Range("H122").Value = 3950
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 3950
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 11850
$ws.Range("N122").Value = -16750
$ws.Range("M122").ClearContents()

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 27963.166
$ws.Range("J59").Value = 27963.166
$ws.Range("L59").Value = 27963.166
$ws.Range("N59").Value = -29657.166
$ws.Range("H86").Value = 1990.8148
$ws.Range("I86").Value = 1738.7273
$ws.Range("J86").Value = 3100
$ws.Range("K86").Value = 1738.7273
$ws.Range("L86").Value = 3100
$ws.Range("M86").Value = -615.7273
$ws.Range("N86").Value = -5346
$ws.Range("H89").Value = 1990.8148
$ws.Range("I89").Value = 1738.7273
$ws.Range("J89").Value = 3100
$ws.Range("K89").Value = 8693.636500000001
$ws.Range("L89").Value = 15500
$ws.Range("M89").Value = -3077.636500000001
$ws.Range("N89").Value = -26732

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2500
$ws.Range("I16").Value = 3250
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 3250
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -2963
$ws.Range("N16").Value = -1574
$ws.Range("H86").Value = 2257.9412
$ws.Range("I86").Value = 1714.625
$ws.Range("K86").Value = 1714.625
$ws.Range("M86").Value = -591.625
$ws.Range("H89").Value = 2257.9412
$ws.Range("I89").Value = 1714.625
$ws.Range("K89").Value = 8573.125
$ws.Range("M89").Value = -2957.125
$ws.Range("H97").Value = 13500
$ws.Range("J97").Value = 13500
$ws.Range("L97").Value = 13500
$ws.Range("N97").Value = -15482
$ws.Range("H99").Value = 1117.5
$ws.Range("I99").Value = 1171
$ws.Range("J99").Value = 957
$ws.Range("K99").Value = 1171
$ws.Range("L99").Value = 957
$ws.Range("M99").Value = 327
$ws.Range("N99").Value = -3953
$ws.Range("H113").Value = 2500
$ws.Range("I113").Value = 3250
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 3250
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = -1080
$ws.Range("N113").Value = -5340
$ws.Range("H126").Value = 1117.5
$ws.Range("I126").Value = 1171
$ws.Range("J126").Value = 957
$ws.Range("K126").Value = 3513
$ws.Range("L126").Value = 2871
$ws.Range("M126").Value = -1043
$ws.Range("N126").Value = -7811

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H24").Value = 826.6667
$ws.Range("I24").Value = 200
$ws.Range("J24").Value = 1140
$ws.Range("K24").Value = 600
$ws.Range("L24").Value = 3420
$ws.Range("M24").Value = -370
$ws.Range("N24").Value = -3880
$ws.Range("H96").Value = 141411410
$ws.Range("J96").Value = 141411410
$ws.Range("L96").Value = 424234230
$ws.Range("N96").Value = -424238348
$ws.Range("H110").Value = 2316.6667
$ws.Range("I110").Value = 2316.6667
$ws.Range("K110").Value = 6950.000100000001
$ws.Range("M110").Value = -2860.000100000001
$ws.Range("H131").Value = 1615345.6
$ws.Range("J131").Value = 1787481.6
$ws.Range("L131").Value = 5362444.800000001
$ws.Range("N131").Value = -5372524.800000001
$ws.Range("H133").Value = 6686.923
$ws.Range("J133").Value = 7799.4736
$ws.Range("L133").Value = 23398.4208
$ws.Range("N133").Value = -33518.4208
$ws.Range("H137").Value = 28048.658
$ws.Range("I137").Value = 4807.5
$ws.Range("J137").Value = 30372.775
$ws.Range("K137").Value = 14422.5
$ws.Range("L137").Value = 91118.32500000001
$ws.Range("M137").Value = -9322.5
$ws.Range("N137").Value = -101318.325

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 502541.1
$ws.Range("J21").Value = 2001331
$ws.Range("L21").Value = 2001331
$ws.Range("N21").Value = -2001677
$ws.Range("H30").Value = 502541.1
$ws.Range("J30").Value = 2001331
$ws.Range("L30").Value = 2001331
$ws.Range("N30").Value = -2001541
$ws.Range("H70").Value = 5194.952
$ws.Range("I70").Value = 5274.0625
$ws.Range("J70").Value = 4941.8
$ws.Range("K70").Value = 5274.0625
$ws.Range("L70").Value = 4941.8
$ws.Range("M70").Value = -5004.0625
$ws.Range("N70").Value = -5481.8
$ws.Range("H73").Value = 5194.952
$ws.Range("I73").Value = 5274.0625
$ws.Range("J73").Value = 4941.8
$ws.Range("K73").Value = 5274.0625
$ws.Range("L73").Value = 4941.8
$ws.Range("M73").Value = -4338.0625
$ws.Range("N73").Value = -6813.8
$ws.Range("H102").Value = 1477.9412
$ws.Range("I102").Value = 1200.8462
$ws.Range("J102").Value = 2378.5
$ws.Range("K102").Value = 1200.8462
$ws.Range("L102").Value = 2378.5
$ws.Range("M102").Value = 421.1538
$ws.Range("N102").Value = -5622.5

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 7073.12
$ws.Range("I122").Value = 8212.895
$ws.Range("J122").Value = 3463.8333
$ws.Range("K122").Value = 24638.685
$ws.Range("L122").Value = 10391.4999
$ws.Range("M122").Value = -22188.685
$ws.Range("N122").Value = -15291.4999

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5354.2856
$ws.Range("J81").Value = 6796
$ws.Range("L81").Value = 13592
$ws.Range("N81").Value = -15714
$ws.Range("H84").Value = 5354.2856
$ws.Range("J84").Value = 6796
$ws.Range("L84").Value = 67960
$ws.Range("N84").Value = -78568
$ws.Range("H126").Value = 3534.2
$ws.Range("I126").Value = 3650.5715
$ws.Range("K126").Value = 10951.7145
$ws.Range("M126").Value = -8481.7145
